# Revise config file handling
# Append the next day's logged config record (row 47) to each of the four
# DE_* sheets, mirroring the existing row layout:
#   A: timestamp (date/time, same NumberFormat as the row above)
#   B-E: hex byte strings (stored as text)
#   F-I: decimal values

$wb = $excel.ActiveWorkbook

function Add-ConfigRow {
    param($SheetName, $TimeVal, $BVal, $CVal, $DVal, $EVal, $FVal, $GVal, $HVal, $IVal)

    $ws = $wb.Worksheets.Item($SheetName)
    $newRow = 47

    $ws.Cells.Item($newRow, 1).Value = $TimeVal
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $BVal
    $ws.Cells.Item($newRow, 3).Value = $CVal
    $ws.Cells.Item($newRow, 4).Value = $DVal
    $ws.Cells.Item($newRow, 5).Value = $EVal

    $ws.Cells.Item($newRow, 6).Value = $FVal
    $ws.Cells.Item($newRow, 7).Value = $GVal
    $ws.Cells.Item($newRow, 8).Value = $HVal
    $ws.Cells.Item($newRow, 9).Value = $IVal
}

$timeVal = 45833.43362268519

# DE_LFT_#1
$g1 = [double]"7.598631275147109e+23"
Add-ConfigRow "DE_LFT_#1" $timeVal "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x64" "0x14" 380 $g1 356 14

# DE_LFT_#2
$g2 = [double]"5.68432987514711e+23"
Add-ConfigRow "DE_LFT_#2" $timeVal "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x68" "0xe" 380 $g2 360 14

# DE_PLT_#1
$g3 = [double]"5.68631262647114e+23"
Add-ConfigRow "DE_PLT_#1" $timeVal "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x7F" "0x7" 130 $g3 127 7

# DE_PLT_#2
$g4 = [double]"9.85046333984776e+23"
Add-ConfigRow "DE_PLT_#2" $timeVal "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x7E" "0x3" 130 $g4 126 3

Write-Output "Row 47 appended to all four sheets"
